$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = 0.4790531655081721
$ws.Range("D5").Value = 0.0563591959421379
$ws.Range("D9").Value = 0.8359947398083788
$ws.Range("D11").Value = 0.1315047905316551
$ws.Range("D13").Value = 0.159684388502724
$ws.Range("D16").Value = 0.8172083411609995
$ws.Range("D19").Value = 0.6481307533345857
$ws.Range("D20").Value = 0.9956791283111027
$ws.Range("D21").Value = 0.1033251925605861
$ws.Range("D25").Value = 3.503663347736239
$ws.Range("D30").Value = 1.437159496524516
$ws.Range("D34").Value = 1.33383430396393
$ws.Range("D36").Value = 1.005072327634792
$ws.Range("D40").Value = 0.8547811384557581
$ws.Range("D41").Value = 1.85985346609055
$ws.Range("D43").Value = 0.02817959797106895
$ws.Range("D51").Value = 0.3381551756528274
$ws.Range("D57").Value = 0.0375727972947586
$ws.Range("D59").Value = 0.3851211722712756
$ws.Range("D61").Value = 0.2536163817396205
$ws.Range("D67").Value = 0.1127183918842758
$ws.Range("D72").Value = 0.09393199323689649
$ws.Range("D73").Value = 0.4884463648318617
$ws.Range("D75").Value = 0.4226939695660342
$ws.Range("D78").Value = 0.1690775878264137
$ws.Range("D88").Value = 0.2630095810633102
$ws.Range("D91").Value = 0.1315047905316551
$ws.Range("D94").Value = 2.066503851211723
$ws.Range("D98").Value = 0.2160435844448619
$ws.Range("D101").Value = 0.6011647567161376
$ws.Range("D107").Value = 2.902498591020101
$ws.Range("D108").Value = 3.015216982904377
$ws.Range("D109").Value = 0.4320871688897238
$ws.Range("D110").Value = 0.3945143715949653
$ws.Range("D112").Value = 1.465339094495585
$ws.Range("D113").Value = 2.846139395077964
$ws.Range("D117").Value = 2.113469847830171
$ws.Range("D118").Value = 4.705992861168514
$ws.Range("D120").Value = 2.658275408604171
$ws.Range("D121").Value = 0.06575239526582755
$ws.Range("D137").Value = 0.06575239526582755
$ws.Range("D140").Value = 0.319368777005448
$ws.Range("D147").Value = 0.375727972947586
$ws.Range("D149").Value = 0.2817959797106894
$ws.Range("D155").Value = 0.8078151418373097
$ws.Range("D159").Value = 0.09393199323689649
$ws.Range("D165").Value = 0.2724027803869998
$ws.Range("D168").Value = 0.187863986473793
$ws.Range("D169").Value = 0.638737554010896
$ws.Range("D170").Value = 0.06575239526582755
$ws.Range("D171").Value = 0.4226939695660342
$ws.Range("D178").Value = 0.9111403343978959
$ws.Range("D182").Value = 0.6011647567161376
$ws.Range("D183").Value = 0.4602667668607928
$ws.Range("D184").Value = 0.1502911891790344
$ws.Range("D190").Value = 0.319368777005448
$ws.Range("D194").Value = 0.8172083411609995
$ws.Range("D196").Value = 1.315047905316551
$ws.Range("D204").Value = 0.4320871688897238
$ws.Range("D215").Value = 0.4039075709186549
$ws.Range("D217").Value = 0.4133007702423445
$ws.Range("D223").Value = 0.2911891790343791
$ws.Range("D226").Value = 1.531091489761413
$ws.Range("D228").Value = 0.3851211722712756
$ws.Range("D233").Value = 0.0187863986473793
$ws.Range("D234").Value = 0.02817959797106895
$ws.Range("D235").Value = 0.0187863986473793
$ws.Range("D238").Value = 0.09393199323689649

# Delete the row for the ADP/SCONJ bigram (now zero-count), shifting subsequent rows up
$ws.Rows.Item(242).Delete()

